# Burndownchart project 2 - "Final submit met alle nieuwe documenten"
#
# The "sprint 3" sheet's "real line" row (row 3) only had the first
# (total points) and last (day5) values filled in; the daily progress
# values for day1..day4 were still missing. Fill them in, matching the
# corresponding "ideal line" columns that are already complete.
#
# Also restore the selected cell / active view for each of the three
# sprint sheets to match where the author last left the cursor.

$wb = $excel.ActiveWorkbook

# --- sprint 1: selection stays on C3 -------------------------------------
$ws1 = $wb.Worksheets.Item("sprint 1")
$ws1.Activate()
$ws1.Range("C3").Select()

# --- sprint 2: selection moves from H3 to G3 (and the H1 scroll/freeze
#     position is cleared so the view resets to the top-left) -------------
$ws2 = $wb.Worksheets.Item("sprint 2")
$ws2.Activate()
$ws2.Range("G3").Select()

# --- sprint 3: fill in the missing "real line" burndown values and move
#     the selection from R12 to D2 -----------------------------------------
$ws3 = $wb.Worksheets.Item("sprint 3")
$ws3.Activate()
$ws3.Range("D3").Value = 43
$ws3.Range("E3").Value = 22
$ws3.Range("F3").Value = 0
$ws3.Range("G3").Value = 0
$ws3.Range("D2").Select()
